$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (col 9), shifting existing columns I:L to J:M.
# This mirrors right-click > Insert on column I in Excel, which also makes
# the new column inherit the width/formatting of the column to its left.
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").ColumnWidth = $ws.Range("H1").ColumnWidth

# New header / type / mapping for the inserted "vacancy_link" column.
$ws.Range("I18").Value = "vacancy_link"
$ws.Range("I19").Value = "varchar"
$ws.Range("I20").Value = "item['alternate_url']"

# Row 20's row-level format flag is dropped on save in the target file while
# every individual cell keeps its style, so stash the row's style, clear the
# row-level flag, then restore the same style cell-by-cell.
$ws.Range("B20").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Rows(20).ClearFormats()
$ws.Range("Z1").Copy()
$ws.Range("B20:M20").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# Match the author's final selection position recorded in the saved file.
$ws.Range("E32").Select()
